$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new year column R (2021), mirroring the existing Q (2020) column.
# Copy the formatting from column Q into column R first.
$ws.Range("Q3:Q8").Copy() | Out-Null
$ws.Range("R3:R8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Header year
$ws.Range("R3").Value = 2021

# Branches per 100,000 adults (formula row)
$ws.Range("R4").Formula = "=R6/R8*100000"

# ATMs per 100,000 adults (formula row)
$ws.Range("R5").Formula = "=R7/R8*100000"

# Total branches of commercial banks
$ws.Range("R6").Value = 312

# Total ATMs
$ws.Range("R7").Value = 1910

# Adult resident population
$ws.Range("R8").Value = 4409166

# Update the view: move the active selection to R15 (also resets the
# scrolled top-left cell back to the default).
$ws.Range("R15").Select() | Out-Null

$wb.Save()
